$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh / re-sort the "Estado de Cuenta" detail table (rows 16-25): the
# database behind it was regenerated, so every worker/period combination now
# carries the unified "Salario Basico" (828116) and the rows are reordered by
# Periodo Mora (ascending: 1912, 2001, 2002, 2003, 2004) with the two workers
# interleaved for each period.

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1049941850"
$ws.Range("D16").Value = "YOHELIS TATIANA ZABALETA GARCIA"
$ws.Range("E16").Value = "1912"
$ws.Range("F16").Value = 33125
$ws.Range("G16").Value = 828116

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1048941394"
$ws.Range("D17").Value = "WENDY MARIA ATENCIO JULIO"
$ws.Range("E17").Value = "1912"
$ws.Range("F17").Value = 33125
$ws.Range("G17").Value = 828116

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1049941850"
$ws.Range("D18").Value = "YOHELIS TATIANA ZABALETA GARCIA"
$ws.Range("E18").Value = "2001"
$ws.Range("F18").Value = 33125
$ws.Range("G18").Value = 828116

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1048941394"
$ws.Range("D19").Value = "WENDY MARIA ATENCIO JULIO"
$ws.Range("E19").Value = "2001"
$ws.Range("F19").Value = 33125
$ws.Range("G19").Value = 828116

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1049941850"
$ws.Range("D20").Value = "YOHELIS TATIANA ZABALETA GARCIA"
$ws.Range("E20").Value = "2002"
$ws.Range("F20").Value = 33125
$ws.Range("G20").Value = 828116

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1048941394"
$ws.Range("D21").Value = "WENDY MARIA ATENCIO JULIO"
$ws.Range("E21").Value = "2002"
$ws.Range("F21").Value = 33125
$ws.Range("G21").Value = 828116

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1049941850"
$ws.Range("D22").Value = "YOHELIS TATIANA ZABALETA GARCIA"
$ws.Range("E22").Value = "2003"
$ws.Range("F22").Value = 33125
$ws.Range("G22").Value = 828116

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1048941394"
$ws.Range("D23").Value = "WENDY MARIA ATENCIO JULIO"
$ws.Range("E23").Value = "2003"
$ws.Range("F23").Value = 33125
$ws.Range("G23").Value = 828116

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1049941850"
$ws.Range("D24").Value = "YOHELIS TATIANA ZABALETA GARCIA"
$ws.Range("E24").Value = "2004"
$ws.Range("F24").Value = 16562
$ws.Range("G24").Value = 828116

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1048941394"
$ws.Range("D25").Value = "WENDY MARIA ATENCIO JULIO"
$ws.Range("E25").Value = "2004"
$ws.Range("F25").Value = 16562
$ws.Range("G25").Value = 828116
